$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 12 new rows above row 1 (shifts all existing content to row+12) ---
$ws.Rows("1:12").Insert()

# --- 2. Fix up conditional formatting ranges (Insert() in this host does not
#        reflow existing conditional-format sqrefs automatically, unlike Excel).
#        ModifyAppliesToRange() on the *existing* rule objects preserves the
#        iconSet/showValue/cfvo settings correctly. ---
$ws.Range("B4:C10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("B16:C22"))
$ws.Range("B12:C12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("B24:C24"))
$ws.Range("H6:H7").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H18:H19"))
$ws.Range("H4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H16"))
$ws.Range("H5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H17"))

# The first rule originally also covered B13:C26 (a second area in the same
# sqref). Re-create that second area (now B25:C38) as its own icon-set rule
# so both ranges keep the conditional formatting.
$extra = $ws.Range("B25:C38").FormatConditions.AddIconSetCondition()
$extra.IconSet = "3Symbols"

# --- 3. Fix up hyperlinks (same reflow issue as conditional formatting) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E18"), "mailto:Customer@beancafe.com")
$ws.Hyperlinks.Add($ws.Range("E16"), "mailto:lm@g.com")
$ws.Hyperlinks.Add($ws.Range("E17"), "mailto:lm@g.com")
$ws.Hyperlinks.Add($ws.Range("E19"), "mailto:Customer@beancafe.com")

# --- 4. Populate the new top rows with the new scenario text. Cell-write
#        order matters because it determines shared-string insertion order. ---
$ws.Range("A1").Value = "If user registers, then person is automatically created."
$ws.Range("A11").Value = "If a person has made a booking before and then registers? "
$ws.Range("B2").Value = "does user email exist in Person table? "
$ws.Range("C3").Value = "if yes, then update user ID"
$ws.Range("C4").Value = "if no then, create person"
$ws.Range("A6").Value = "If user "

# --- 5. Column A width + selected cell ---
$ws.Columns("A").ColumnWidth = 4.6
$ws.Range("C4").Select()
